# Glottolog update (v4.8 -> v5.0): refresh the cached "datetimeFigureOut"
# placeholder text on every master/layout/notes-master, and bump the six
# cross-table counts that shifted as a result of the re-tabulation.

function Get-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -eq $name) { return $sh }
    }
    return $null
}

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$p = $ppt.ActivePresentation
$newDate = "3/14/24"

# --- 1. Slide master date placeholder -------------------------------------
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# --- 2. Every slide layout's date placeholder ------------------------------
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# --- 3. Notes master date placeholder --------------------------------------
Set-DatePlaceholderText $p.NotesMaster.Shapes $newDate

# --- 4. The six cross-table counts on slide 1 ------------------------------
$slide = $p.Slides.Item(1)
$group84 = Get-ShapeByName $slide.Shapes "Group 84"

$tableUpdates = @(
    @{ Shapes = $slide.Shapes;          Name = "Table 1";  Row = 2; Col = 2; Text = "2973" }
    @{ Shapes = $group84.GroupItems;    Name = "Table 4";  Row = 2; Col = 2; Text = "2990" }
    @{ Shapes = $group84.GroupItems;    Name = "Table 19"; Row = 2; Col = 1; Text = "3253" }
    @{ Shapes = $group84.GroupItems;    Name = "Table 21"; Row = 2; Col = 1; Text = "3227" }
    @{ Shapes = $slide.Shapes;          Name = "Table 78"; Row = 2; Col = 1; Text = "3134" }
    @{ Shapes = $slide.Shapes;          Name = "Table 80"; Row = 2; Col = 1; Text = "3222" }
)

foreach ($u in $tableUpdates) {
    $tblShape = Get-ShapeByName $u.Shapes $u.Name
    $tblShape.Table.Cell($u.Row, $u.Col).Shape.TextFrame.TextRange.Text = $u.Text
}

Write-Output "Applied glottolog v4.8 -> v5.0 updates"
